$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append after the current last row (row 95)
$newRows = @(
    @{ Row = 96; Values = @(44336, 905, 2611, 650, 5977, 450, 2209, 1000, 28500, 0, 0, 0, 3674, 509, 1140, 42970) },
    @{ Row = 97; Values = @(44337, 350, 2511, 755, 5827, 517, 2002, 0, 28500, 0, 0, 0, 3674, 24, 1153, 42514) }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $vals = $entry.Values
    for ($col = 1; $col -le $vals.Count; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.Value = $vals[$col - 1]
    }
    # Column A uses the date-formatted style (same as existing rows, numFmt yyyy-mm-dd)
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd"
}
